$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$apos = [char]0x2019

# Row 3 - Astute Financial Management
$ws.Range("B3").Value = "Operating since 2000, Astute has established itself as a highly respected financial services provider with a membership base comprising of approximately 500 finance, insurance and financial planning specialists across Australia.`nWe are a privately owned business, leading the market with a fully integrated financial services offering to our members. With a loan book of over `$20 billion dollars, we are well positioned for growth."
$ws.Range("C3").Value = "Toowong"
$ws.Range("D3").Value = 27.4843
$ws.Range("E3").Value = 152.9837
$ws.Range("F3").Value = 1000000
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "1000000                                                                                                                              "
$ws.Range("G3").NumberFormat = "#,##0"

# Row 4 - Aussie
$ws.Range("B4").Value = "Aussie works closely with industry and government to ensure a strong customer first industry. Our team are passionate about continuing to find ways to improve the home loan experience. We${apos}ll give you the confidence to make a move."
$ws.Range("C4").Value = "philippines"
$ws.Range("D4").Value = 12.8797
$ws.Range("E4").Value = 121.774
$ws.Range("F4").Value = 135000
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "2100000                                                                                                                              "
$ws.Range("G4").NumberFormat = "#,##0"
